# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamp cells with fresh report-generation
# timestamps (stored as text, matching the existing text-formatted values).

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the 5705e82c... row (G3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-20 18:55:18"

# zh-cn sheet: "Correspond Handoff Datetime" (H3) and
# "Correspond Handback DateTime" (K3) for the 5705e82c... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-20 18:55:14"
$wsZhCn.Range("K3").Value = "2016-08-20 18:55:30"

# de-de sheet: "Correspond Handoff Datetime" (H3) and
# "Correspond Handback DateTime" (K3) for the 5705e82c... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-20 18:55:18"
$wsDeDe.Range("K3").Value = "2016-08-20 18:55:36"
